{"js": "// Replace the underlined section label \"Laws:\" with \"Law:\" (singular).\nconst results = context.document.body.search(\"Laws:\", { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(\"Law:\", Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Replace the underlined section label \"Laws:\" with \"Law:\" (singular).\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n\n# Find.Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,\n#              MatchAllWordForms, Forward, Wrap, Format, ReplaceWith, Replace)\n$find.Execute(\"Laws:\", $true, $false, $false, $false, $false, $true, 1, $false, \"Law:\", 2) | Out-Null\n"}
